$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "struggle"
$ws.Cells.Item(2, 3).Value = 2.571562051773072
$ws.Cells.Item(2, 4).Value = -5.523353099822998
$ws.Cells.Item(2, 5).Value = -4.892651081085205
$ws.Cells.Item(2, 6).Value = 0.0096211275085806
$ws.Cells.Item(2, 7).Value = 0.0177150927484035
$ws.Cells.Item(2, 8).Value = -0.0091629782691597

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "struggle"
$ws.Cells.Item(3, 3).Value = 2.563363254070282
$ws.Cells.Item(3, 4).Value = -5.564052700996399
$ws.Cells.Item(3, 5).Value = -4.925167679786682
$ws.Cells.Item(3, 6).Value = -0.0056505035609006
$ws.Cells.Item(3, 7).Value = -0.007177666760981
$ws.Cells.Item(3, 8).Value = 0

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "struggle"
$ws.Cells.Item(4, 3).Value = 2.571200489997864
$ws.Cells.Item(4, 4).Value = -5.45090651512146
$ws.Cells.Item(4, 5).Value = -4.94497549533844
$ws.Cells.Item(4, 6).Value = -0.0161879286170005
$ws.Cells.Item(4, 7).Value = 0.0122173046693205
$ws.Cells.Item(4, 8).Value = -0.0047342055477201

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "struggle"
$ws.Cells.Item(5, 3).Value = 2.582025349140167
$ws.Cells.Item(5, 4).Value = -5.429405391216278
$ws.Cells.Item(5, 5).Value = -4.891633093357086
$ws.Cells.Item(5, 6).Value = 0.0029016099870204
$ws.Cells.Item(5, 7).Value = -0.0010690141934901
$ws.Cells.Item(5, 8).Value = -0.009468411095440299

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "struggle"
$ws.Cells.Item(6, 3).Value = 2.521161556243896
$ws.Cells.Item(6, 4).Value = -5.436496257781982
$ws.Cells.Item(6, 5).Value = -4.74793529510498
$ws.Cells.Item(6, 6).Value = 0.00534507073462
$ws.Cells.Item(6, 7).Value = 0.0088575463742017
$ws.Cells.Item(6, 8).Value = 0.0045814891345798

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "struggle"
$ws.Cells.Item(7, 3).Value = 2.451793313026428
$ws.Cells.Item(7, 4).Value = -5.718138635158539
$ws.Cells.Item(7, 5).Value = -4.555151760578156
$ws.Cells.Item(7, 6).Value = 0.0036651915870606
$ws.Cells.Item(7, 7).Value = -0.0503963828086853
$ws.Cells.Item(7, 8).Value = 0.0389426611363887

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "struggle"
$ws.Cells.Item(8, 3).Value = 4.312116026878357
$ws.Cells.Item(8, 4).Value = -5.115874052047729
$ws.Cells.Item(8, 5).Value = -7.913362622261047
$ws.Cells.Item(8, 6).Value = -1.588860511779785
$ws.Cells.Item(8, 7).Value = -2.905885934829712
$ws.Cells.Item(8, 8).Value = 0.6783658266067505

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "struggle"
$ws.Cells.Item(9, 3).Value = 3.11747863888741
$ws.Cells.Item(9, 4).Value = -5.436717748641967
$ws.Cells.Item(9, 5).Value = -5.15502381324769
$ws.Cells.Item(9, 6).Value = -0.9424123764038086
$ws.Cells.Item(9, 7).Value = -3.084411382675171
$ws.Cells.Item(9, 8).Value = 0.4729624092578888

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "struggle"
$ws.Cells.Item(10, 3).Value = 3.307996869087219
$ws.Cells.Item(10, 4).Value = -5.476139068603516
$ws.Cells.Item(10, 5).Value = -0.757482767105103
$ws.Cells.Item(10, 6).Value = 0.2102903574705124
$ws.Cells.Item(10, 7).Value = -1.237612962722778
$ws.Cells.Item(10, 8).Value = 0.2101376503705978

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "struggle"
$ws.Cells.Item(11, 3).Value = 4.280841529369356
$ws.Cells.Item(11, 4).Value = -4.743481069803241
$ws.Cells.Item(11, 5).Value = 2.116865754127498
$ws.Cells.Item(11, 6).Value = -0.4138612151145935
$ws.Cells.Item(11, 7).Value = -0.9240863919258118
$ws.Cells.Item(11, 8).Value = 0.252134621143341

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "struggle"
$ws.Cells.Item(12, 3).Value = 1.148419260978671
$ws.Cells.Item(12, 4).Value = -4.571187555789954
$ws.Cells.Item(12, 5).Value = 3.019693136215206
$ws.Cells.Item(12, 6).Value = 0.432645320892334
$ws.Cells.Item(12, 7).Value = 1.750434398651123
$ws.Cells.Item(12, 8).Value = -0.1994474977254867

$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "struggle"
$ws.Cells.Item(13, 3).Value = -5.23933637142183
$ws.Cells.Item(13, 4).Value = -7.058503150939933
$ws.Cells.Item(13, 5).Value = -0.1314393877982891
$ws.Cells.Item(13, 6).Value = 1.822669148445129
$ws.Cells.Item(13, 7).Value = 5.886908531188965
$ws.Cells.Item(13, 8).Value = -2.078774452209473

$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "struggle"
$ws.Cells.Item(14, 3).Value = 5.029654502868652
$ws.Cells.Item(14, 4).Value = -12.30031204223633
$ws.Cells.Item(14, 5).Value = -12.83510589599609
$ws.Cells.Item(14, 6).Value = 0.7883216142654419
$ws.Cells.Item(14, 7).Value = 2.360841512680054
$ws.Cells.Item(14, 8).Value = 1.160491228103638

$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "struggle"
$ws.Cells.Item(15, 3).Value = -1.211387872695933
$ws.Cells.Item(15, 4).Value = -5.499476730823506
$ws.Cells.Item(15, 5).Value = -6.533813059329978
$ws.Cells.Item(15, 6).Value = 0.5829181671142578
$ws.Cells.Item(15, 7).Value = 1.261283993721008
$ws.Cells.Item(15, 8).Value = 0.6953173875808716

$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "struggle"
$ws.Cells.Item(16, 3).Value = -3.101414084434511
$ws.Cells.Item(16, 4).Value = -3.493300497531889
$ws.Cells.Item(16, 5).Value = -7.406978726387027
$ws.Cells.Item(16, 6).Value = -0.3068070709705353
$ws.Cells.Item(16, 7).Value = 1.792431354522705
$ws.Cells.Item(16, 8).Value = -0.2423607856035232

$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "struggle"
$ws.Cells.Item(17, 3).Value = -3.497859179973603
$ws.Cells.Item(17, 4).Value = -3.212260961532593
$ws.Cells.Item(17, 5).Value = -10.54216539859774
$ws.Cells.Item(17, 6).Value = 0.2667953968048095
$ws.Cells.Item(17, 7).Value = -0.7637342810630798
$ws.Cells.Item(17, 8).Value = 0.0548251569271087

$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "struggle"
$ws.Cells.Item(18, 3).Value = -1.64927089214325
$ws.Cells.Item(18, 4).Value = -7.43125307559967
$ws.Cells.Item(18, 5).Value = -8.658325910568239
$ws.Cells.Item(18, 6).Value = -1.207069754600525
$ws.Cells.Item(18, 7).Value = -6.127589225769043
$ws.Cells.Item(18, 8).Value = -1.29213273525238

$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "struggle"
$ws.Cells.Item(19, 3).Value = -1.025534451007847
$ws.Cells.Item(19, 4).Value = -11.15312564373016
$ws.Cells.Item(19, 5).Value = -2.677394092082968
$ws.Cells.Item(19, 6).Value = -1.138041973114014
$ws.Cells.Item(19, 7).Value = -3.366478443145752
$ws.Cells.Item(19, 8).Value = 0.9393580555915833

$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "struggle"
$ws.Cells.Item(20, 3).Value = -1.654657959938065
$ws.Cells.Item(20, 4).Value = -9.341210365295419
$ws.Cells.Item(20, 5).Value = 0.4500467777252322
$ws.Cells.Item(20, 6).Value = 1.038929104804993
$ws.Cells.Item(20, 7).Value = -0.7336491346359253
$ws.Cells.Item(20, 8).Value = -1.004567861557007

$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "struggle"
$ws.Cells.Item(21, 3).Value = 2.828480809926987
$ws.Cells.Item(21, 4).Value = -5.927687406539919
$ws.Cells.Item(21, 5).Value = -2.06907200813294
$ws.Cells.Item(21, 6).Value = -0.5413793325424194
$ws.Cells.Item(21, 7).Value = -1.692554831504822
$ws.Cells.Item(21, 8).Value = -0.2292271852493286

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "struggle"
$ws.Cells.Item(22, 3).Value = 3.216025352478028
$ws.Cells.Item(22, 4).Value = -4.024631500244141
$ws.Cells.Item(22, 5).Value = 4.593602180480957
$ws.Cells.Item(22, 6).Value = -0.3520110845565796
$ws.Cells.Item(22, 7).Value = 0.6711881756782532
$ws.Cells.Item(22, 8).Value = -0.102472648024559

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "struggle"
$ws.Cells.Item(23, 3).Value = -1.550865292549141
$ws.Cells.Item(23, 4).Value = -5.700099587440493
$ws.Cells.Item(23, 5).Value = 3.356234908103941
$ws.Cells.Item(23, 6).Value = 2.986520290374756
$ws.Cells.Item(23, 7).Value = 4.694652080535889
$ws.Cells.Item(23, 8).Value = 0.0166460778564214

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "struggle"
$ws.Cells.Item(24, 3).Value = -2.655673027038576
$ws.Cells.Item(24, 4).Value = -3.905611395835868
$ws.Cells.Item(24, 5).Value = -3.318085908889795
$ws.Cells.Item(24, 6).Value = 1.599092483520508
$ws.Cells.Item(24, 7).Value = 4.835456371307373
$ws.Cells.Item(24, 8).Value = 1.327257513999939

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "struggle"
$ws.Cells.Item(25, 3).Value = -1.430848956108096
$ws.Cells.Item(25, 4).Value = -1.578429281711582
$ws.Cells.Item(25, 5).Value = -8.265119194984431
$ws.Cells.Item(25, 6).Value = -0.3246748745441437
$ws.Cells.Item(25, 7).Value = 0.3605632185935974
$ws.Cells.Item(25, 8).Value = 0.3962988257408142

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "struggle"
$ws.Cells.Item(26, 3).Value = -2.407764434814463
$ws.Cells.Item(26, 4).Value = -1.766093373298651
$ws.Cells.Item(26, 5).Value = -8.389460563659659
$ws.Cells.Item(26, 6).Value = 0.05283984541893
$ws.Cells.Item(26, 7).Value = 0.2293798923492431
$ws.Cells.Item(26, 8).Value = -0.2738203406333923

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "struggle"
$ws.Cells.Item(27, 3).Value = -2.611050009727472
$ws.Cells.Item(27, 4).Value = -2.573673054575919
$ws.Cells.Item(27, 5).Value = -8.294337868690496
$ws.Cells.Item(27, 6).Value = -0.6475171446800232
$ws.Cells.Item(27, 7).Value = -0.2813034355640411
$ws.Cells.Item(27, 8).Value = 0.0232128798961639

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "struggle"
$ws.Cells.Item(28, 3).Value = 0.61596310138703
$ws.Cells.Item(28, 4).Value = -2.870795279741297
$ws.Cells.Item(28, 5).Value = -8.901223957538587
$ws.Cells.Item(28, 6).Value = -0.1372919678688049
$ws.Cells.Item(28, 7).Value = -2.705674886703491
$ws.Cells.Item(28, 8).Value = -0.5198463201522827

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "struggle"
$ws.Cells.Item(29, 3).Value = 2.079445004463198
$ws.Cells.Item(29, 4).Value = -5.368536770343783
$ws.Cells.Item(29, 5).Value = -3.524431616067886
$ws.Cells.Item(29, 6).Value = -0.6478226184844971
$ws.Cells.Item(29, 7).Value = -0.2128865420818328
$ws.Cells.Item(29, 8).Value = -0.0656680166721344

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "struggle"
$ws.Cells.Item(30, 3).Value = 0.3216586112976074
$ws.Cells.Item(30, 4).Value = -3.676267147064209
$ws.Cells.Item(30, 5).Value = -3.865855693817138
$ws.Cells.Item(30, 6).Value = -0.1000291854143142
$ws.Cells.Item(30, 7).Value = 0.1372919678688049
$ws.Cells.Item(30, 8).Value = -0.1838704347610473

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "struggle"
$ws.Cells.Item(31, 3).Value = 1.706132471561434
$ws.Cells.Item(31, 4).Value = -4.47040206193924
$ws.Cells.Item(31, 5).Value = -5.197765350341799
$ws.Cells.Item(31, 6).Value = 0.2654209434986114
$ws.Cells.Item(31, 7).Value = 0.0520762614905834
$ws.Cells.Item(31, 8).Value = 0.0438295826315879

